# BUG: Fix read_excel w/parse_cols & empty dataset (#23661)
# Adds a third worksheet ("Sheet3") containing only a single header row
# (A, B, C, D, E, F) and no data rows, so read_excel can be exercised
# against an otherwise-empty dataset.

$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet so it lands at the end
# of the tab strip (Sheet1, Sheet2, Sheet3).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Header-only row; shared strings A-D already exist in the workbook, E/F
# are new.
$ws3.Range("A1").Value = "A"
$ws3.Range("B1").Value = "B"
$ws3.Range("C1").Value = "C"
$ws3.Range("D1").Value = "D"
$ws3.Range("E1").Value = "E"
$ws3.Range("F1").Value = "F"

# Leave the cursor parked one row below the last header cell, matching
# the saved selection state (F2) on the new active sheet.
[void]$ws3.Range("F2").Select()
